# Update gh-pages to output generated at 456a3b4
# Sheet "展览" (index 1): F3 226->228, F4 2544->2550, F5 40->41
# Sheet "全部类型" (index 4): F5 226->228, F6 2544->2550, F7 40->41

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item(1)
$wsExhibition.Range("F3").Value = 228
$wsExhibition.Range("F4").Value = 2550
$wsExhibition.Range("F5").Value = 41

$wsAllTypes = $wb.Worksheets.Item(4)
$wsAllTypes.Range("F5").Value = 228
$wsAllTypes.Range("F6").Value = 2550
$wsAllTypes.Range("F7").Value = 41
